# Updates cryptos list prices (column D) and 1h volume % (column E)
# to reflect the latest scraped values, per commit
# "Updated cryptos list on Wed Sep 27 20:13:18 UTC 2023 with GitHub Actions"
#
# Column D ("Price") values are plain text (e.g. "26.250.26") rather than
# numbers, so a leading apostrophe is used to force Excel to keep them as
# text, and the quote-prefix style that introduces is immediately reset
# back to Normal so no stray cell formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.250.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "'1.597.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'211.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").Value = "'0.504"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").Value = "'18.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.45%  "
$ws.Range("D12").Value = "'1.822.28"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").Value = "'1.596.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.30%  "
$ws.Range("D14").Value = "'4.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("E15").Value = "  -2.40%  "
$ws.Range("D16").Value = "'63.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("D17").Value = "'26.262.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "'229.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.08%  "
$ws.Range("E19").Value = "  +3.05%  "
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("E22").Value = "  -0.63%  "
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("E25").Value = "  +1.03%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").Value = "'15.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.23%  "
$ws.Range("D30").Value = "'0.0492"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.55%  "
$ws.Range("E31").Value = "  -0.29%  "
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("D33").Value = "'1.469.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.64%  "
$ws.Range("D34").Value = "'2.94"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  -0.33%  "
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("D37").Value = "'0.566"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.33%  "
$ws.Range("E38").Value = "  -1.07%  "
$ws.Range("D39").Value = "'0.820"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.33%  "
$ws.Range("E40").Value = "  -2.70%  "
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("D42").Value = "'2.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.91%  "
$ws.Range("D43").Value = "'0.934"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.03%  "
$ws.Range("D44").Value = "'1.734.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.69%  "
$ws.Range("D45").Value = "'0.758"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.89%  "
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("D47").Value = "'87.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("E48").Value = "  -1.08%  "
$ws.Range("D49").Value = "'0.0501"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("E51").Value = "  +0.23%  "
